$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the brand-new strings first, in the same order they were
# originally authored, so the shared-strings table comes out in the
# same sequence as the source workbook.
$ws.Range("F9").Value  = "Restricted access to productive and financial resources"
$ws.Range("F15").Value = "?"
$ws.Range("F4").Value  = "Not present. Generated it?"
$ws.Range("F5").Value  = "ok (removed Code in the title)"
$ws.Range("F11").Value = "?. Not divided"
$ws.Range("F14").Value = "Secure access to land assets - law + Secure access to land assets - practice?"

# Remaining column F cells, reusing the strings above (and the
# existing "ok" string already present in the workbook).
$ws.Range("F3").Value  = "ok"
$ws.Range("F6").Value  = "Not present. Generated it?"
$ws.Range("F7").Value  = "ok"
$ws.Range("F8").Value  = "Not present. Generated it?"
$ws.Range("F10").Value = "Not present. Generated it?"
$ws.Range("F12").Value = "?. Not divided"
$ws.Range("F13").Value = "?. Not divided"

# New header for column F (row 2) - plain year number, bold + centered
# like the rest of the header row.
$ws.Range("F2").Value = 2019
$ws.Range("F2").Font.Bold = $true
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4108

# Body cells of column F: centered, matching the rest of the table.
$bodyRange = $ws.Range("F3:F15")
$bodyRange.HorizontalAlignment = -4108
$bodyRange.VerticalAlignment = -4108

# Column F width, widened to fit the new text
$ws.Columns.Item(6).ColumnWidth = 48.917

# Update selection to match the new active cell in the authored file
$ws.Range("F14").Select()
